$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 59458.668
$ws.Range("I11").Value = 59458.668
$ws.Range("K11").Value = 59458.668
$ws.Range("M11").Value = -59318.668

$ws.Range("H17").Value = 1785.8846
$ws.Range("J17").Value = 1785.8846
$ws.Range("L17").Value = 5357.6538
$ws.Range("N17").Value = -5693.6538

$ws.Range("H51").Value = 1499.8889
$ws.Range("I51").Value = 1499
$ws.Range("K51").Value = 1499
$ws.Range("M51").Value = -1015

$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()

$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()

$ws.Range("H111").Value = 624.75
$ws.Range("I111").Value = 500
$ws.Range("K111").Value = 1500
$ws.Range("M111").Value = 1567

$ws.Range("H141").Value = 6666.3335
$ws.Range("I141").Value = 6666.3335
$ws.Range("K141").Value = 19999.0005
$ws.Range("M141").Value = -14819.0005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 15464.25
$ws.Range("I74").Value = 15464.25
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 15464.25
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -14590.25
$ws.Range("N74").ClearContents()

$ws.Range("H77").Value = 15464.25
$ws.Range("I77").Value = 15464.25
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 77321.25
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -72953.25
$ws.Range("N77").ClearContents()

$ws.Range("H88").Value = 3625
$ws.Range("I88").Value = 2000
$ws.Range("J88").Value = 4166.6665
$ws.Range("K88").Value = 2000
$ws.Range("L88").Value = 4166.6665
$ws.Range("M88").Value = -1594
$ws.Range("N88").Value = -4978.6665

$ws.Range("H91").Value = 3625
$ws.Range("I91").Value = 2000
$ws.Range("J91").Value = 4166.6665
$ws.Range("K91").Value = 2000
$ws.Range("L91").Value = 4166.6665
$ws.Range("M91").Value = -596
$ws.Range("N91").Value = -6974.6665

$ws.Range("H122").Value = 3726
$ws.Range("I122").Value = 3726
$ws.Range("K122").Value = 11178
$ws.Range("M122").Value = -8728

$ws.Range("H132").Value = 3993.3333
$ws.Range("I132").Value = 3988.4
$ws.Range("K132").Value = 11965.2
$ws.Range("M132").Value = -9435.200000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 7425.5
$ws.Range("J86").Value = 13153.5
$ws.Range("L86").Value = 13153.5
$ws.Range("N86").Value = -15399.5

$ws.Range("H89").Value = 7425.5
$ws.Range("J89").Value = 13153.5
$ws.Range("L89").Value = 65767.5
$ws.Range("N89").Value = -76999.5

$ws.Range("H134").Value = 7447.273
$ws.Range("I134").Value = 4904.8887
$ws.Range("K134").Value = 14714.6661
$ws.Range("M134").Value = -12179.6661

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 92.111115
$ws.Range("J7").Value = 97.8
$ws.Range("L7").Value = 97.8
$ws.Range("N7").Value = -323.8

$ws.Range("H22").Value = 354.0909
$ws.Range("I22").Value = 369.21054
$ws.Range("J22").Value = 258.33334
$ws.Range("K22").Value = 369.21054
$ws.Range("L22").Value = 258.33334
$ws.Range("M22").Value = -19.21053999999998
$ws.Range("N22").Value = -958.33334

$ws.Range("H62").Value = 17601
$ws.Range("J62").Value = 17601
$ws.Range("L62").Value = 17601
$ws.Range("N62").Value = -18849

$ws.Range("H65").Value = 17601
$ws.Range("J65").Value = 17601
$ws.Range("L65").Value = 88005
$ws.Range("N65").Value = -94245

$ws.Range("H86").Value = 7282.4287
$ws.Range("I86").Value = 7164.5
$ws.Range("K86").Value = 7164.5
$ws.Range("M86").Value = -6041.5

$ws.Range("H89").Value = 7282.4287
$ws.Range("I89").Value = 7164.5
$ws.Range("K89").Value = 35822.5
$ws.Range("M89").Value = -30206.5

$ws.Range("H94").Value = 2697
$ws.Range("I94").Value = 2697
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 2697
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -2246
$ws.Range("N94").ClearContents()

$ws.Range("H99").Value = 3434.4443
$ws.Range("I99").Value = 3285
$ws.Range("K99").Value = 3285
$ws.Range("M99").Value = -1787

$ws.Range("H105").Value = 3402.7778
$ws.Range("I105").Value = 3187.5
$ws.Range("K105").Value = 3187.5
$ws.Range("M105").Value = -1440.5

$ws.Range("H126").Value = 3434.4443
$ws.Range("I126").Value = 3285
$ws.Range("K126").Value = 9855
$ws.Range("M126").Value = -7385

$ws.Range("H132").Value = 4875.778
$ws.Range("I132").Value = 4721
$ws.Range("J132").Value = 4999.6
$ws.Range("K132").Value = 14163
$ws.Range("L132").Value = 14998.8
$ws.Range("M132").Value = -11633
$ws.Range("N132").Value = -20058.8

$ws.Range("H134").Value = 8133
$ws.Range("I134").Value = 8133
$ws.Range("K134").Value = 24399
$ws.Range("M134").Value = -21864

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H16").Value = 100
$ws.Range("I16").Value = 100
$ws.Range("K16").Value = 300
$ws.Range("M16").Value = -127

$ws.Range("H112").Value = 7468.75
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 7468.75
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 22406.25
$ws.Range("N112").Value = -24622.25
$ws.Range("M112").ClearContents()

$ws.Range("H132").Value = 3496.5
$ws.Range("I132").Value = 3595.8
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 32362.2
$ws.Range("L132").Value = 27000
$ws.Range("M132").Value = -29832.2
$ws.Range("N132").Value = -32060

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 5666.6665
$ws.Range("I122").Value = 4500
$ws.Range("J122").Value = 8000
$ws.Range("K122").Value = 13500
$ws.Range("L122").Value = 24000
$ws.Range("M122").Value = -11050
$ws.Range("N122").Value = -28900

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 13224.583
$ws.Range("I22").Value = 17939.2
$ws.Range("J22").Value = 9857
$ws.Range("K22").Value = 17939.2
$ws.Range("L22").Value = 9857
$ws.Range("M22").Value = -17644.2
$ws.Range("N22").Value = -10447

$ws.Range("H27").Value = 13224.583
$ws.Range("I27").Value = 17939.2
$ws.Range("J27").Value = 9857
$ws.Range("K27").Value = 17939.2
$ws.Range("L27").Value = 9857
$ws.Range("M27").Value = -17832.2
$ws.Range("N27").Value = -10071

$ws.Range("H42").Value = 9999
$ws.Range("I42").Value = 9999
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 9999
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = -9436
$ws.Range("N42").ClearContents()

$ws.Range("H49").Value = 9999
$ws.Range("I49").Value = 9999
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 9999
$ws.Range("L49").Value = 0
$ws.Range("M49").Value = -9852
$ws.Range("N49").ClearContents()

$ws.Range("H55").Value = 3112.5
$ws.Range("I55").Value = 3833.3333
$ws.Range("J55").Value = 950
$ws.Range("K55").Value = 3833.3333
$ws.Range("L55").Value = 950
$ws.Range("M55").Value = -3660.3333
$ws.Range("N55").Value = -1296

$ws.Range("H110").Value = 76667
$ws.Range("J110").Value = 76667
$ws.Range("L110").Value = 76667
$ws.Range("N110").Value = -84847

$ws.Range("H136").Value = 1670
$ws.Range("I136").Value = 1670
$ws.Range("K136").Value = 5010
$ws.Range("M136").Value = -2460

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1439.9333
$ws.Range("I100").Value = 1182.125
$ws.Range("K100").Value = 2364.25
$ws.Range("M100").Value = -1823.25

$ws.Range("H122").Value = 2661.6667
$ws.Range("I122").Value = 2661.6667
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 7985.000100000001
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -5535.000100000001
$ws.Range("N122").ClearContents()

$ws.Range("H126").Value = 3367
$ws.Range("I126").Value = 3400.5
$ws.Range("J126").Value = 3300
$ws.Range("K126").Value = 10201.5
$ws.Range("L126").Value = 9900
$ws.Range("M126").Value = -7731.5
$ws.Range("N126").Value = -14840
